$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new log entry in row 17, matching the date of row 16,
# with 2 hours spent and a new description "formatting, descriptions".

# Copy formatting (number format/styles) from row 16 down to row 17 first.
$ws.Range("A16:C16").Copy() | Out-Null
$ws.Range("A17:C17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Use Value2 to get/set the raw numeric date serial (avoids locale/formatted strings).
$ws.Range("A17").Value2 = $ws.Range("A16").Value2
$ws.Range("B17").Value2 = 2
$ws.Range("C17").Value = "formatting, descriptions"

# Update the active selection as in the diff.
$ws.Range("B20").Select() | Out-Null
